$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117, pushing the existing rows 117-223 down to 118-224.
$ws.Rows(117).Insert()

# Populate the newly inserted row 117 with the new record (it was blank after Insert,
# apart from the inherited date style on column D).
$ws.Cells.Item(117, 1).Value = 5
$ws.Cells.Item(117, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(117, 3).Value = "Maule"
$ws.Cells.Item(117, 4).Value = 44512
$ws.Cells.Item(117, 5).Value = 7
$ws.Cells.Item(117, 6).Value = 100112023
$ws.Cells.Item(117, 7).Value = "Brócoli"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 4000
$ws.Cells.Item(117, 11).Value = 600
$ws.Cells.Item(117, 12).Value = 600
$ws.Cells.Item(117, 13).Value = 600
$ws.Cells.Item(117, 14).Value = "`$/unidad"
$ws.Cells.Item(117, 15).Value = "Región del Maule"
$ws.Cells.Item(117, 16).Value = 600
$ws.Cells.Item(117, 17).Value = 1
$ws.Cells.Item(117, 18).Value = "Hortaliza"
